$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G6").Value = "Dharwad"
$ws.Range("G11").Value = "Shivamogga (Shimoga)"
$ws.Range("G12").Value = "Chikkamagaluru (Chikmagalur)"
$ws.Range("G19").Value = "Chikkamagaluru (Chikmagalur)"
$ws.Range("G33").Value = "Shivamogga (Shimoga)"
$ws.Range("G34").Value = "Shivamogga (Shimoga)"
$ws.Range("G35").Value = "Shivamogga (Shimoga)"
$ws.Range("G36").Value = "Shivamogga (Shimoga)"
